function Escape-Xml($text) {
    $text = $text -replace '&', '&amp;'
    $text = $text -replace '<', '&lt;'
    $text = $text -replace '>', '&gt;'
    return $text
}

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the opening
#    "Play Before Time Runs Out for free today!" heading.
# ------------------------------------------------------------------
$metaDesc = ": Read our review of Before Time Runs Out, an immersive slot game set in oriental culture, and play for free with exciting bonus games and free spins."

$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaXml = '<w:p ' + $wns + '><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>' + (Escape-Xml $metaDesc) + '</w:t></w:r></w:p>'
$metaPara.Range.InsertXML($metaXml)

# ------------------------------------------------------------------
# 2) Near the end of the document: drop the duplicated bold
#    "Play Before Time Runs Out for free today!" paragraph, and turn
#    the following italic paragraph into the image-generation prompt.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$playAgainPara = $d.Paragraphs($count - 1)
$playAgainPara.Range.Delete()

$promptText = 'Prompt: Create a feature image for the game "Before Time Runs Out" that portrays a happy Maya warrior with glasses in a cartoon style. To capture the essence of the game, the warrior should be surrounded by elements of oriental culture, such as a scimitar, a palace, and an hourglass. The image should be bright and colorful, with a starry sky in the background and swirls like on the game''s console frame. Ensure the image depicts excitement, adventure, and offers a glimpse of the fantasy world the game creates.'

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$promptXml = '<w:p ' + $wns + '><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>' + (Escape-Xml $promptText) + '</w:t></w:r></w:p>'
$lastPara.Range.InsertXML($promptXml)
